# Updates cryptos list price/volume figures (and a RocketPoolETH/FraxShare
# row-order swap) to match the latest scrape, per GitHub Actions commit.
# Numeric-looking "Price" text values are written with a leading apostrophe
# so Excel keeps them as text (matching the original inlineStr cells)
# instead of auto-converting them to numbers; the style is reset to
# "Normal" right after so no stray number-format/quote-prefix style is
# left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.075.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.57%  "

$ws.Range("D3").Value = "'1.791.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.71%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'223.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.34%  "

$ws.Range("E6").Value = "  -0.67%  "

$ws.Range("D8").Value = "'32.27"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.89%  "

$ws.Range("D9").Value = "'0.284"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.20%  "

$ws.Range("D10").Value = "'0.0708"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.46%  "

$ws.Range("E11").Value = "  +0.01%  "

$ws.Range("D12").Value = "'2.049.03"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.72%  "

$ws.Range("D13").Value = "'1.790.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.39%  "

$ws.Range("D14").Value = "'10.82"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.65%  "

$ws.Range("D15").Value = "'0.624"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.98%  "

$ws.Range("D16").Value = "'34.033.06"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.80%  "

$ws.Range("D17").Value = "'4.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.91%  "

$ws.Range("D18").Value = "'67.92"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.91%  "

$ws.Range("D19").Value = "'242.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.94%  "

$ws.Range("D20").Value = "'0.0₃0781"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.46%  "

$ws.Range("E21").Value = "  +0.18%  "

$ws.Range("D22").Value = "'10.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.02%  "

$ws.Range("D23").Value = "'4.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.32%  "

$ws.Range("E24").Value = "  -2.80%  "

$ws.Range("D25").Value = "'158.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.76%  "

$ws.Range("D26").Value = "'16.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.01%  "

$ws.Range("D27").Value = "'6.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.26%  "

$ws.Range("E28").Value = "  -2.13%  "

$ws.Range("E29").Value = "  +0.13%  "

$ws.Range("D30").Value = "'0.0516"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.31%  "

$ws.Range("E31").Value = "  +0.19%  "

$ws.Range("D32").Value = "'3.64"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.20%  "

$ws.Range("E33").Value = "  -3.85%  "

$ws.Range("E34").Value = "  -4.35%  "

$ws.Range("D35").Value = "'1.385.37"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.37%  "

$ws.Range("D36").Value = "'0.645"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.26%  "

$ws.Range("E37").Value = "  -1.91%  "

$ws.Range("E38").Value = "  -3.89%  "

$ws.Range("E39").Value = "  -0.01%  "

$ws.Range("E40").Value = "  -6.65%  "

$ws.Range("E41").Value = "  -3.66%  "

$ws.Range("D42").Value = "'0.912"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.77%  "

$ws.Range("D43").Value = "'2.16"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.10%  "

$ws.Range("D44").Value = "'0.0₆0137"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.37%  "

$ws.Range("E45").Value = "  -0.27%  "

$ws.Range("E46").Value = "  -0.70%  "

$ws.Range("D47").Value = "'107.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.64%  "

$ws.Range("E50").Value = "  +0.01%  "

$ws.Range("D51").Value = "'11.90"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.05%  "

# Row 48/49 content swap (RocketPoolETH <-> FraxShare) with updated values
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "'5.82"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.53%  "

$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "'1.947.71"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.50%  "
